$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -3.789114356040955
$ws.Cells.Item(2, 4).Value = 5.144325375556946
$ws.Cells.Item(2, 5).Value = -3.264913499355316
$ws.Cells.Item(2, 6).Value = 0.0216857157647609
$ws.Cells.Item(2, 7).Value = -0.0339030213654041
$ws.Cells.Item(2, 8).Value = 0.0201585534960031

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -3.717108857631684
$ws.Cells.Item(3, 4).Value = 5.147888684272766
$ws.Cells.Item(3, 5).Value = -3.339606630802154
$ws.Cells.Item(3, 6).Value = -0.0120645882561802
$ws.Cells.Item(3, 7).Value = -0.0719293802976608
$ws.Cells.Item(3, 8).Value = 0.0242818929255008

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -3.921339607238769
$ws.Cells.Item(4, 4).Value = 5.102599048614501
$ws.Cells.Item(4, 5).Value = -3.249480080604553
$ws.Cells.Item(4, 6).Value = 0.0125227374956011
$ws.Cells.Item(4, 7).Value = -0.0106901414692401
$ws.Cells.Item(4, 8).Value = -0.0226020142436027

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -3.555192089080811
$ws.Cells.Item(5, 4).Value = 4.907798504829406
$ws.Cells.Item(5, 5).Value = -2.923101136088372
$ws.Cells.Item(5, 6).Value = -0.0018325957935303
$ws.Cells.Item(5, 7).Value = -0.0296269636601209
$ws.Cells.Item(5, 8).Value = -0.0087048299610614

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -3.523949909210205
$ws.Cells.Item(6, 4).Value = 4.786228704452514
$ws.Cells.Item(6, 5).Value = -1.791911107301711
$ws.Cells.Item(6, 6).Value = 0.0019853119738399
$ws.Cells.Item(6, 7).Value = 0.0355829000473022
$ws.Cells.Item(6, 8).Value = 0.0200058370828628

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -4.734984159469604
$ws.Cells.Item(7, 4).Value = 4.268091917037964
$ws.Cells.Item(7, 5).Value = -0.5654808729887018
$ws.Cells.Item(7, 6).Value = 0.0106901414692401
$ws.Cells.Item(7, 7).Value = 0.08765916526317589
$ws.Cells.Item(7, 8).Value = 0.0525344125926494

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = -5.463105487823487
$ws.Cells.Item(8, 4).Value = 3.789239883422852
$ws.Cells.Item(8, 5).Value = 0.125382423400879
$ws.Cells.Item(8, 6).Value = -0.0532979927957057
$ws.Cells.Item(8, 7).Value = 0.3320052623748779
$ws.Cells.Item(8, 8).Value = 0.1557706445455551

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -6.153488552570344
$ws.Cells.Item(9, 4).Value = 2.863918662071227
$ws.Cells.Item(9, 5).Value = -0.2511623546481139
$ws.Cells.Item(9, 6).Value = -0.1586722433567047
$ws.Cells.Item(9, 7).Value = 0.3882048726081848
$ws.Cells.Item(9, 8).Value = -0.0940732508897781

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -7.365730333328246
$ws.Cells.Item(10, 4).Value = 1.174339866638185
$ws.Cells.Item(10, 5).Value = 0.3271868914365764
$ws.Cells.Item(10, 6).Value = -0.2698497176170349
$ws.Cells.Item(10, 7).Value = 0.6501133441925049
$ws.Cells.Item(10, 8).Value = -0.4689917862415313

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "falling"
$ws.Cells.Item(11, 3).Value = -7.554701018333435
$ws.Cells.Item(11, 4).Value = 1.320136770606041
$ws.Cells.Item(11, 5).Value = 0.9918349981307982
$ws.Cells.Item(11, 6).Value = 0.0048869219608604
$ws.Cells.Item(11, 7).Value = 0.8135197758674622
$ws.Cells.Item(11, 8).Value = -1.29075825214386

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "falling"
$ws.Cells.Item(12, 3).Value = -6.541906356811527
$ws.Cells.Item(12, 4).Value = 3.541345179080956
$ws.Cells.Item(12, 5).Value = -0.3953665494918761
$ws.Cells.Item(12, 6).Value = 0.8275696635246277
$ws.Cells.Item(12, 7).Value = -2.131614208221436
$ws.Cells.Item(12, 8).Value = -2.630385637283325

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "falling"
$ws.Cells.Item(13, 3).Value = 41.03047697544085
$ws.Cells.Item(13, 4).Value = 6.894529008865351
$ws.Cells.Item(13, 5).Value = -6.942315888404829
$ws.Cells.Item(13, 6).Value = 1.503492116928101
$ws.Cells.Item(13, 7).Value = -4.445571899414063
$ws.Cells.Item(13, 8).Value = -2.242486238479614

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "falling"
$ws.Cells.Item(14, 3).Value = 64.96824674606313
$ws.Cells.Item(14, 4).Value = 6.919658172130577
$ws.Cells.Item(14, 5).Value = -9.210757869482023
$ws.Cells.Item(14, 6).Value = -1.950034618377685
$ws.Cells.Item(14, 7).Value = 3.770718336105346
$ws.Cells.Item(14, 8).Value = 0.2194533348083496

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "falling"
$ws.Cells.Item(15, 3).Value = -6.831269806623445
$ws.Cells.Item(15, 4).Value = 3.269432669878014
$ws.Cells.Item(15, 5).Value = -0.1760272979736253
$ws.Cells.Item(15, 6).Value = 6.251442432403564
$ws.Cells.Item(15, 7).Value = -2.222022294998169
$ws.Cells.Item(15, 8).Value = -3.614794969558716

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "falling"
$ws.Cells.Item(16, 3).Value = -2.175523471832268
$ws.Cells.Item(16, 4).Value = 5.449279594421388
$ws.Cells.Item(16, 5).Value = 1.586600971221922
$ws.Cells.Item(16, 6).Value = 0.8747590184211731
$ws.Cells.Item(16, 7).Value = 2.045024156570435
$ws.Cells.Item(16, 8).Value = 1.580461144447327

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "falling"
$ws.Cells.Item(17, 3).Value = 2.27016156911849
$ws.Cells.Item(17, 4).Value = 5.673981070518497
$ws.Cells.Item(17, 5).Value = 0.4227319359779442
$ws.Cells.Item(17, 6).Value = -0.0520762614905834
$ws.Cells.Item(17, 7).Value = -0.5412266254425049
$ws.Cells.Item(17, 8).Value = 0.422871470451355

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "falling"
$ws.Cells.Item(18, 3).Value = 4.514346599578862
$ws.Cells.Item(18, 4).Value = 5.118150138854979
$ws.Cells.Item(18, 5).Value = -1.915704894065858
$ws.Cells.Item(18, 6).Value = -0.2570215463638305
$ws.Cells.Item(18, 7).Value = -7.184691429138184
$ws.Cells.Item(18, 8).Value = 2.630538463592529

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "falling"
$ws.Cells.Item(19, 3).Value = -6.036725759506252
$ws.Cells.Item(19, 4).Value = 9.008556652069101
$ws.Cells.Item(19, 5).Value = -0.9792453408241246
$ws.Cells.Item(19, 6).Value = 2.87503719329834
$ws.Cells.Item(19, 7).Value = -1.23394775390625
$ws.Cells.Item(19, 8).Value = -0.1574505120515823

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "falling"
$ws.Cells.Item(20, 3).Value = -1.178059291839593
$ws.Cells.Item(20, 4).Value = 4.555799674987787
$ws.Cells.Item(20, 5).Value = 3.518247509002693
$ws.Cells.Item(20, 6).Value = -0.157145082950592
$ws.Cells.Item(20, 7).Value = 1.678810358047485
$ws.Cells.Item(20, 8).Value = 1.309847831726074

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "falling"
$ws.Cells.Item(21, 3).Value = -1.642422831058502
$ws.Cells.Item(21, 4).Value = 6.837077736854553
$ws.Cells.Item(21, 5).Value = 3.844039088487625
$ws.Cells.Item(21, 6).Value = 0.982424020767212
$ws.Cells.Item(21, 7).Value = 1.2322678565979
$ws.Cells.Item(21, 8).Value = -0.3689626157283783

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "falling"
$ws.Cells.Item(22, 3).Value = 0.4585734605789122
$ws.Cells.Item(22, 4).Value = 6.029543757438664
$ws.Cells.Item(22, 5).Value = 3.976011931896209
$ws.Cells.Item(22, 6).Value = 0.0326812900602817
$ws.Cells.Item(22, 7).Value = 0.007177666760981
$ws.Cells.Item(22, 8).Value = -0.1241583600640297

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "falling"
$ws.Cells.Item(23, 3).Value = 0.8665444850921629
$ws.Cells.Item(23, 4).Value = 5.183062970638275
$ws.Cells.Item(23, 5).Value = 4.418135178089142
$ws.Cells.Item(23, 6).Value = -0.00167987938039
$ws.Cells.Item(23, 7).Value = 0.0308486949652433
$ws.Cells.Item(23, 8).Value = 0.0305432621389627

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "falling"
$ws.Cells.Item(24, 3).Value = 0.2307996749877924
$ws.Cells.Item(24, 4).Value = 4.777379417419433
$ws.Cells.Item(24, 5).Value = 4.424502086639404
$ws.Cells.Item(24, 6).Value = 0.012980886735022
$ws.Cells.Item(24, 7).Value = -0.06704246252775189
$ws.Cells.Item(24, 8).Value = -0.0529925599694252

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "falling"
$ws.Cells.Item(25, 3).Value = -0.31689715385437
$ws.Cells.Item(25, 4).Value = 4.877218794822693
$ws.Cells.Item(25, 5).Value = 4.225887775421143
$ws.Cells.Item(25, 6).Value = 0.0004581489483825
$ws.Cells.Item(25, 7).Value = -0.0694859251379966
$ws.Cells.Item(25, 8).Value = -0.005192354787141

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "falling"
$ws.Cells.Item(26, 3).Value = -0.4965919017791747
$ws.Cells.Item(26, 4).Value = 5.156597185134888
$ws.Cells.Item(26, 5).Value = 4.267510080337524
$ws.Cells.Item(26, 6).Value = -0.015118914656341
$ws.Cells.Item(26, 7).Value = -0.0235183127224445
$ws.Cells.Item(26, 8).Value = 0.047036625444889

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "falling"
$ws.Cells.Item(27, 3).Value = -0.5277259349822996
$ws.Cells.Item(27, 4).Value = 5.028044939041139
$ws.Cells.Item(27, 5).Value = 4.214276224374771
$ws.Cells.Item(27, 6).Value = 0.0102319931611418
$ws.Cells.Item(27, 7).Value = 0.0080939643085002
$ws.Cells.Item(27, 8).Value = -0.0006108652451075

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "falling"
$ws.Cells.Item(28, 3).Value = -0.6842149734497071
$ws.Cells.Item(28, 4).Value = 4.842281723022461
$ws.Cells.Item(28, 5).Value = 3.947901606559754
$ws.Cells.Item(28, 6).Value = 0.0326812900602817
$ws.Cells.Item(28, 7).Value = 0.0154243474826216
$ws.Cells.Item(28, 8).Value = 0.0415388382971286

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "falling"
$ws.Cells.Item(29, 3).Value = -0.5955237627029424
$ws.Cells.Item(29, 4).Value = 5.176781976222992
$ws.Cells.Item(29, 5).Value = 4.094217467308044
$ws.Cells.Item(29, 6).Value = 0.0032070425804704
$ws.Cells.Item(29, 7).Value = 0.0175623763352632
$ws.Cells.Item(29, 8).Value = -0.0198531206697225

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "falling"
$ws.Cells.Item(30, 3).Value = -0.1085210800170896
$ws.Cells.Item(30, 4).Value = 5.183276605606079
$ws.Cells.Item(30, 5).Value = 4.1697988986969
$ws.Cells.Item(30, 6).Value = 0.0097738439217209
$ws.Cells.Item(30, 7).Value = -0.0273362193256616
$ws.Cells.Item(30, 8).Value = -0.0233655963093042

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = -0.2790900230407717
$ws.Cells.Item(31, 4).Value = 5.053058326244355
$ws.Cells.Item(31, 5).Value = 4.114621889591216
$ws.Cells.Item(31, 6).Value = 0.0073303831741213
$ws.Cells.Item(31, 7).Value = 0.0030543261673301
$ws.Cells.Item(31, 8).Value = -0.0368046313524246
